# Update the "Estado de Cuenta" worksheet with the new EC database values
# and trim it down to a single page of data (part 1 of the new EC batch).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Trim the worker-detail table down to 3 rows -------------------------
# The old table had 13 detail rows (16-28) before 4 blank rows and the
# signature footer (33-34). The new extract only keeps 3 detail rows
# (16-18), so delete the old rows 19-28; this leaves the blank-row gap
# (4 rows) intact and slides the footer up from 33/34 to 23/24 automatically,
# preserving its original formatting/merged cells.
$ws.Range("B19:J28").EntireRow.Delete() | Out-Null

# --- Header / summary block ---------------------------------------------
$ws.Range("E11").Value = 67353

$ws.Range("C13").Value = 2
$ws.Range("F13").Value = 2

# --- Worker detail table (rows 16-18) -----------------------------------
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1047459519"
$ws.Range("D16").Value = "CAROLINA SANTAMARIA MOLINA"
$ws.Range("E16").Value = "1907"
$ws.Range("F16").Value = 30916
$ws.Range("G16").Value = 908526

$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1047459519"
$ws.Range("D17").Value = "CAROLINA SANTAMARIA MOLINA"
$ws.Range("E17").Value = "1908"
$ws.Range("F17").Value = 33125
$ws.Range("G17").Value = 908526

$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1143398464"
$ws.Range("D18").Value = "LAURA DANIELA ROMERO LEON"
$ws.Range("E18").Value = "1908"
$ws.Range("F18").Value = 3312
$ws.Range("G18").Value = 877803

# --- Footer (signature block), now directly below the data table --------
$ws.Range("B23").Value = "___________________________________"
$ws.Range("H23").Value = "___________________________________"
$ws.Range("B24").Value = "NOMBRE DEL REPRESENTANTE LEGAL"
$ws.Range("H24").Value = "FIRMA DEL REPRESENTANTE LEGAL"
